$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 423 (shifts existing rows 423:440 down to 425:442)
$ws.Range("A423:A424").EntireRow.Insert()

# --- New row 423 (Coliflor, Primera, week of 2021-11-09) ---
$ws.Range("A423").Value = 8
$ws.Range("B423").Value = "Terminal La Palmera de La Serena"
$ws.Range("C423").Value = "Coquimbo"
$ws.Range("D423").Value = 44509
$ws.Range("E423").Value = 4
$ws.Range("F423").Value = 100112008
$ws.Range("G423").Value = "Coliflor"
$ws.Range("H423").Value = "Sin especificar"
$ws.Range("I423").Value = "Primera"
$ws.Range("J423").Value = 2200
$ws.Range("K423").Value = 600
$ws.Range("L423").Value = 700
$ws.Range("M423").Value = 650
$ws.Range("N423").Value = "$/unidad"
$ws.Range("O423").Value = "Provincia del Elquí"
$ws.Range("P423").Value = 650
$ws.Range("Q423").Value = 1
$ws.Range("R423").Value = "Hortaliza"

# --- New row 424 (Coliflor, Segunda, week of 2021-11-09) ---
$ws.Range("A424").Value = 8
$ws.Range("B424").Value = "Terminal La Palmera de La Serena"
$ws.Range("C424").Value = "Coquimbo"
$ws.Range("D424").Value = 44509
$ws.Range("E424").Value = 4
$ws.Range("F424").Value = 100112008
$ws.Range("G424").Value = "Coliflor"
$ws.Range("H424").Value = "Sin especificar"
$ws.Range("I424").Value = "Segunda"
$ws.Range("J424").Value = 1460
$ws.Range("K424").Value = 500
$ws.Range("L424").Value = 550
$ws.Range("M424").Value = 525
$ws.Range("N424").Value = "$/unidad"
$ws.Range("O424").Value = "Provincia del Elquí"
$ws.Range("P424").Value = 525
$ws.Range("Q424").Value = 1
$ws.Range("R424").Value = "Hortaliza"
